$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings
# (e.g. "1.000", "0.4943") are stored verbatim as text, not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.693.31'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '1.893.37'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = '311.22'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").Value = '0.4943'
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("D8").Value = '0.3787'
$ws.Range("E8").Value = '  -1.47%  '
$ws.Range("D9").Value = '0.07315'
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").Value = '0.9073'
$ws.Range("E10").Value = '  -5.02%  '
$ws.Range("D11").Value = '20.58'
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").Value = '0.07639'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").Value = '1.870.10'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").Value = '5.463'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '6.635'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '91.05'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '0.000008728'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '27.694.13'
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = '14.44'
$ws.Range("E21").Value = '  -4.33%  '
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("D23").Value = '2.112.90'
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '154.15'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").Value = '1.846'
$ws.Range("E26").Value = '  -5.24%  '
$ws.Range("D27").Value = '18.38'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").Value = '2.154'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = '115.30'
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("E30").Value = '  -3.71%  '
$ws.Range("D31").Value = '0.08943'
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").Value = '3.200'
$ws.Range("E32").Value = '  -4.19%  '
$ws.Range("E33").Value = '  -2.55%  '
$ws.Range("D34").Value = '0.7636'
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("D35").Value = '4.625'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("D36").Value = '2.558'
$ws.Range("E36").Value = '  -7.86%  '
$ws.Range("D37").Value = '0.02039'
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("D39").Value = '0.05279'
$ws.Range("E39").Value = '  -2.36%  '
$ws.Range("D40").Value = '0.5469'
$ws.Range("E40").Value = '  -2.52%  '
$ws.Range("D41").Value = '2.988'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").Value = '6.862'
$ws.Range("E42").Value = '  -4.06%  '
$ws.Range("D43").Value = '8.518'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("D44").Value = '112.48'
$ws.Range("E44").Value = '  +4.22%  '
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").Value = '10.55'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").Value = '1.634'
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("D50").Value = '67.13'
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("D51").Value = '0.06055'
$ws.Range("E51").Value = '  -1.41%  '

# Restore column D to its original (default) style so no stray
# number-format styling is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
